$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (style matches existing header cells -> reuse H1's format)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-16 for columns I (I0) and J (IF)
$data = @(
    @(2, 1, 4),
    @(3, 1, 6),
    @(4, 1, 7),
    @(5, 1, 5),
    @(6, 1, 7),
    @(7, 1, 6),
    @(8, 1, 6),
    @(9, 1, 7),
    @(10, 1, 6),
    @(11, 1, 7),
    @(12, 1, 7),
    @(13, 1, 6),
    @(14, 1, 5),
    @(15, 3, 3),
    @(16, 1, 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $if = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}
